$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the refreshed Price (column D) and Volume(1h) (column E) values.
# Each value is written with a leading apostrophe so Excel stores it as
# literal text (matching the original inlineStr text cells) rather than
# converting numeric-looking strings (e.g. "328.31") or percentages
# (e.g. "5.91%") into real numbers. The style is reset right after each
# write so the cell keeps the workbook default formatting instead of the
# text quote-prefix style Excel applies automatically for such entries.

$ws.Range("D2").Value = "'328.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'5.91%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'40.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'7.67%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.662"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'10.53%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08108"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.38%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.550"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'3.45%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'8.689"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'1.960"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'4.29%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E10").Value = "'2.78%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1272"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'8.33%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1994"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'5.20%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09178"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'3.40%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.03551"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'7.21%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.09610"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.02%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001310"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-4.80%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.006074"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.97%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.371"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.75%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3506"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.38%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.499"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'17.61%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1401"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'8.47%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2506"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'4.09%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04432"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.98%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001253"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'4.39%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004338"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.24%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-14.97%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003993"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'37.65%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02527"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'16.97%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05212"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'4.09%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007789"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'2.69%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1431"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'5.50%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.009059"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'6.78%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'9.06%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.01051"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'31.62%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006731"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'2.40%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.08%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002874"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-12.75%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'59.32%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.08%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.08%"
$ws.Range("E51").Style = "Normal"
